$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking text values in column D need a quote-prefix trick so Excel
# keeps storing them as text (matching the original inlineStr cells) instead of
# silently converting to a number; the cell style is saved/restored so no style
# index churn is introduced.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.Value = "'" + $value
    $rng.Style = $origStyle
}

Set-TextValue "D2" "269.16"
Set-TextValue "D3" "21.12"
Set-TextValue "D4" "6.252"
Set-TextValue "D5" "0.06203"
Set-TextValue "D6" "3.565"
Set-TextValue "D7" "6.538"
Set-TextValue "D8" "1.440"
Set-TextValue "D9" "0.8240"
Set-TextValue "D10" "0.1656"
Set-TextValue "D11" "0.08255"
Set-TextValue "D12" "0.03558"
Set-TextValue "D13" "0.03186"
Set-TextValue "D14" "0.09190"
Set-TextValue "D15" "3.764"
Set-TextValue "D16" "0.001627"
Set-TextValue "D17" "0.04680"
Set-TextValue "D18" "0.006338"
Set-TextValue "D19" "0.006196"
Set-TextValue "D20" "0.001068"
Set-TextValue "D22" "3.728"
Set-TextValue "D23" "2.261"
Set-TextValue "D24" "0.01377"
Set-TextValue "D28" "0.0002714"
Set-TextValue "D40" "0.04701"
Set-TextValue "D41" "0.007012"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.004201"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1120"
$ws.Range("E43").Value = "42BKEXTokenBKK"
Set-TextValue "D44" "0.01127"
Set-TextValue "D45" "0.00006227"
Set-TextValue "D46" "0.0009903"
Set-TextValue "D47" "0.00000000750"
Set-TextValue "D48" "0.9907"
Set-TextValue "D50" "0.00001900"
Set-TextValue "D51" "0.01240"
